$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing strings (Float -> Normal, Boost -> Hizli) ---
$ws.Range("A17").Value2 = "Normal Şarj"
$ws.Range("A18").Value2 = "Hızlı Şarj"
$ws.Range("E25").Value2 = "Normal ve hızlı testleri tamamlanmış varsayılıyor."
$ws.Range("B29").Value2 = "Normal ve hızlı testleri tamamlanmış varsayılıyor."
$ws.Range("E26").Value2 = "Manuel olarak normal şarj moduna geç."
$ws.Range("B30").Value2 = "Manuel olarak normal şarj moduna geç."
$ws.Range("B34").Value2 = "Hızlı şarj(Zmn) iletisini oku."
$ws.Range("B36").Value2 = "Yazılı kayıtlardaki hızlı şarj değerlerine uygun olduğuna bak."

# --- Add new rows 37-39 with new content ---
$ws.Range("B37").Value2 = "Sürenin dolmasını bekle."
$ws.Range("B38").Value2 = "Normal şarj(Zmn) iletisini oku."
$ws.Range("B39").Value2 = "Yazılı kayıtlardaki normal şarj değerlerine uygun olduğuna bak."

# --- Update selection to B41 ---
$ws.Range("B41").Select()
